$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.95
$ws.Range("H2").Value = 2.92
$ws.Range("I2").Value = 4.3

$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 2.18
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 7.7
$ws.Range("V2").Value = 9.25
$ws.Range("W2").Value = 17
$ws.Range("Y2").Value = 45
$ws.Range("Z2").Value = 5.9
$ws.Range("AA2").Value = 6
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 150
$ws.Range("AE2").Value = 22
$ws.Range("AG2").Value = 80
